# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# Add a new player "G.Gilbert" to the QB sheet (row 5), with all stat columns
# (B:L) set to 0 -- same shape as the existing rows.
$qb = $wb.Worksheets.Item("QB")
$qb.Range("A5").Value = "G.Gilbert"
$qb.Range("B5:L5").Value = 0

# Make QB the active/selected sheet (it becomes tabSelected) and select
# the newly-added stat range.
$qb.Activate()
$qb.Range("B4:L5").Select()

# The RB sheet was previously the tab-selected one; it no longer is after
# QB becomes active (Excel keeps RB's own selection at K6 untouched).
